# Updates the odds values on the active worksheet (sheet1) to reflect
# the latest FlashScore odds snapshot for 2024-11-25.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Central Cordoba vs Rosario Central)
$ws.Range("G2").Value = 2.35
$ws.Range("I2").Value = 3.6
$ws.Range("W2").Value = 5.5
$ws.Range("X2").Value = 9.5
$ws.Range("Z2").Value = 23
$ws.Range("AE2").Value = 21
$ws.Range("AI2").Value = 15
$ws.Range("AW2").Value = 5

# Row 3 (2 de Mayo vs Sp. Luqueno)
$ws.Range("U3").Value = 2.2
$ws.Range("V3").Value = 1.62

# Row 4 (Ameliano vs General Caballero JLM)
$ws.Range("G4").Value = 2.38
$ws.Range("I4").Value = 2.88
$ws.Range("J4").Value = 3.25
$ws.Range("L4").Value = 3.75
$ws.Range("U4").Value = 2.1
$ws.Range("V4").Value = 1.67
$ws.Range("W4").Value = 6.5
$ws.Range("Y4").Value = 10
$ws.Range("Z4").Value = 23
$ws.Range("AJ4").Value = 12
$ws.Range("AK4").Value = 34
$ws.Range("AL4").Value = 29
$ws.Range("AN4").Value = 4.33
$ws.Range("AX4").Value = 19
$ws.Range("AY4").Value = 34
$ws.Range("AZ4").Value = 67
$ws.Range("BB4").Value = 301

# Row 5 (Libertad Asuncion vs Sol de America)
$ws.Range("G5").Value = 2.35
$ws.Range("J5").Value = 3.2
$ws.Range("M5").Value = 1.07
$ws.Range("N5").Value = 9
$ws.Range("O5").Value = 1.36
$ws.Range("P5").Value = 3
$ws.Range("Q5").Value = 2.1
$ws.Range("R5").Value = 1.7
$ws.Range("Y5").Value = 10
$ws.Range("AB5").Value = 34
$ws.Range("AC5").Value = 9
$ws.Range("AE5").Value = 17
$ws.Range("AH5").Value = 8
$ws.Range("AN5").Value = 4.5
$ws.Range("AO5").Value = 15
$ws.Range("AP5").Value = 26
$ws.Range("AQ5").Value = 51
$ws.Range("AY5").Value = 29
